$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.771.37'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.475.56'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = "'321.13"
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").Value = "'92.43"
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = "'0.508"
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("D10").Value = "'32.97"
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").Value = "'0.0856"
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("D13").Value = '2.856.99'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").Value = '2.496.77'
$ws.Range("E16").Value = '  -1.05%  '
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").Value = '41.719.21'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = "'6.46"
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '0.0₃0941'
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("D21").Value = "'71.81"
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("D22").Value = "'11.23"
$ws.Range("E22").Value = '  -3.06%  '
$ws.Range("D23").Value = "'239.80"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  +1.31%  '
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("E28").Value = '  -1.70%  '
$ws.Range("D29").Value = "'9.72"
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("D30").Value = "'36.25"
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("D31").Value = "'155.28"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").Value = "'5.43"
$ws.Range("E32").Value = '  -1.43%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = "'0.0765"
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").Value = "'2.56"
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("D36").Value = "'17.06"
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = "'2.92"
$ws.Range("E37").Value = '  +1.09%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = "'1.84"
$ws.Range("E38").Value = '  +2.33%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("E41").Value = '  -0.71%  '
$ws.Range("D42").Value = "'2.36"
$ws.Range("E42").Value = '  -6.22%  '
$ws.Range("D43").Value = '2.005.11'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("D44").Value = "'0.0282"
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").Value = "'18.71"
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("D47").Value = "'9.44"
$ws.Range("E47").Value = '  +4.28%  '
$ws.Range("D48").Value = '2.735.32'
$ws.Range("E48").Value = '  +1.26%  '
$ws.Range("D49").Value = "'97.49"
$ws.Range("E49").Value = '  +0.36%  '
$ws.Range("D50").Value = "'75.92"
$ws.Range("E50").Value = '  +4.06%  '
$ws.Range("D51").Value = "'67.05"
$ws.Range("E51").Value = '  -0.11%  '
